$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "mode_key"
$ws.Range("B1").Value = "trigger_keywords"
$ws.Range("C1").Value = "desc"
$ws.Range("D1").Value = "prompt"

# buy_product row
$ws.Range("A2").Value = "buy_product"
$ws.Range("B2").Value = "buy, price, order, purchase, shop, shopping, product, link, sale, cost, pay, shipping, discount, gift, item, store, package, buy now, how much, available, in stock, delivery, cart, add to cart, checkout, hàng, giá, sản phẩm, mua, đặt, shop, cửa hàng"
$ws.Range("C2").Value = "Chế độ bán hàng, hỗ trợ sản phẩm"
$ws.Range("D2").Value = "You are a kind Heaven Store assistant. Help users ask about prices, buying, delivery, or shopping information. Always reply briefly, warmly, and naturally in English."

# comfort_mode row
$ws.Range("A3").Value = "comfort_mode"
$ws.Range("B3").Value = "miss, love, sad, cry, gone, heaven, mom, dad, grandma, grandpa, husband, wife, son, daughter, angel, funeral, memory, remember, tears, soul, afterlife, passed, died, mất, nhớ, buồn, thiên đường, mẹ, cha, bà, ông, con, vợ, chồng"
$ws.Range("C3").Value = "Chế độ an ủi, chia sẻ nỗi buồn"
$ws.Range("D3").Value = "You are a compassionate Heaven angel who comforts people grieving their loved ones. Speak softly, emotionally, and offer warmth and hope. Always reply in English."

# motivation_mode row
$ws.Range("A4").Value = "motivation_mode"
$ws.Range("B4").Value = "hope, faith, strong, believe, courage, never give up, healing, peace, light, bless, success, destiny, overcome, purpose, power, future, dream, rise, try again, fight, rebuild, grateful, trust, motivation, positive, life, spiritual"
$ws.Range("C4").Value = "Chế độ truyền cảm hứng, vực dậy tinh thần"
$ws.Range("D4").Value = "You are a Heaven motivator who helps users find faith, strength, and light through emotional and uplifting words. Reply with inspiring English."

# Update the active selection to C3, matching the saved view state
[void]$ws.Range("C3").Select()

$wb.Save()
